{"js": "// Apply EMV2Instantiation.docx documentation update:\n// 1. Append explanatory sentences about where propagation path instances\n//    are kept (EMV2AnnexInstance) to the end of the\n//    \"Propagation path instances are created...\" paragraph.\n// 2. Add three new paragraphs under the \"Property Associations\" heading\n//    describing how EMV2 property associations / copies are created.\n\nconst body = context.document.body;\n\n// --- Change 1: append text to the \"Propagation path instances are created...\" paragraph ---\nconst firstParaResults = body.search(\n  \"Propagation path instances are created for connection instances, user defined propagation path instances, and for bindings.\",\n  { matchCase: true }\n);\nfirstParaResults.load(\"text\");\nawait context.sync();\n\nif (firstParaResults.items.length > 0) {\n  const appendText =\n    \" Propagation path instances are kept in the EMV2AnnexInstance of the component that contains the connection instance or user defined propagation path instance. In the case of bindings the propagation path instance is kept in the EMV2AnnexInstance of the system instance.\";\n  firstParaResults.items[0].insertText(appendText, Word.InsertLocation.end);\n  await context.sync();\n}\n\n// --- Change 2: insert the three new \"Property Associations\" paragraphs ---\nconst headingResults = body.search(\"Property Associations\", { matchCase: true });\nawait context.sync();\n\nif (headingResults.items.length > 0) {\n  const headingPara = headingResults.items[0].paragraphs.getFirst();\n  // The heading is followed by three empty \"Bod\" paragraphs in the original\n  // document. The first stays empty; the next two receive the first two new\n  // paragraphs of text, and a new paragraph (plus a fresh trailing empty\n  // paragraph) is inserted for the remaining content so the document still\n  // ends on an empty paragraph, as in the original.\n  const emptyPara1 = headingPara.getNext();\n  const emptyPara2 = emptyPara1.getNext();\n  const emptyPara3 = emptyPara2.getNext();\n\n  emptyPara2.insertText(\n    \"We create a copy of EMV2 properties in the EMV2 instance as owned property association of the respective EMV2 instance object. We do this for error events, error sources, error propagations.\",\n    Word.InsertLocation.replace\n  );\n  emptyPara3.insertText(\n    \"Currently we create copies of Hazards, OccurrenceDistribution, Severity, and Likelihood properties. These are the only ones being used by safety analyses. There is a generic method to create an instance copy of any EMV2 property named as string.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n\n  // Keep a trailing empty paragraph at the end of the document (matching the\n  // original structure) by inserting it first, then placing the final new\n  // paragraph of text right before it.\n  const trailingEmptyPara = emptyPara3.insertParagraph(\"\", Word.InsertLocation.after);\n  await context.sync();\n\n  trailingEmptyPara.insertParagraph(\n    \"Note that property values can be associated with individual error types of the above mentioned EMV2 objects. We represent error events, error sources with multiple error types (listed as type sets) as separate ConstrainedInstanceObject for each type token representing each error type.  The property association is owned by the appropriate ConstrainedInstanceObject.\",\n    Word.InsertLocation.before\n  );\n  await context.sync();\n}\n", "ps1": "# Apply EMV2Instantiation.docx documentation update:\n# 1. Append explanatory sentences about where propagation path instances\n#    are kept (EMV2AnnexInstance) to the end of the\n#    \"Propagation path instances are created...\" paragraph.\n# 2. Add three new paragraphs under the \"Property Associations\" heading\n#    describing how EMV2 property associations / copies are created.\n\nfunction Get-ParaText($para) {\n    # Paragraph Range.Text includes the trailing paragraph mark (CR / cell\n    # mark) character(s); strip them so we can compare against plain text.\n    return $para.Range.Text.TrimEnd([char]13, [char]7)\n}\n\n$d = $word.ActiveDocument\n\n# --- Change 1: append text to the \"Propagation path instances are created...\" paragraph ---\n$paragraphs = $d.Paragraphs\n$count = $paragraphs.Count\n$targetIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $txt = Get-ParaText $paragraphs.Item($i)\n    if ($txt -eq \"Propagation path instances are created for connection instances, user defined propagation path instances, and for bindings.\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -ge 1) {\n    $p = $d.Paragraphs.Item($targetIndex)\n    $p.Range.InsertAfter(\" Propagation path instances are kept in the EMV2AnnexInstance of the component that contains the connection instance or user defined propagation path instance. In the case of bindings the propagation path instance is kept in the EMV2AnnexInstance of the system instance.\")\n}\n\n# --- Change 2: insert the three new \"Property Associations\" paragraphs ---\n$paragraphs2 = $d.Paragraphs\n$count2 = $paragraphs2.Count\n$headingIndex = -1\nfor ($i = 1; $i -le $count2; $i++) {\n    $txt = Get-ParaText $paragraphs2.Item($i)\n    if ($txt -eq \"Property Associations\") {\n        $headingIndex = $i\n        break\n    }\n}\n\nif ($headingIndex -ge 1) {\n    # The heading is followed by three empty \"Bod\" paragraphs in the\n    # original document. The first stays empty; the next two receive the\n    # first two new paragraphs of text.\n    $p1 = $d.Paragraphs.Item($headingIndex + 1)\n    $p2 = $d.Paragraphs.Item($headingIndex + 2)\n    $p3 = $d.Paragraphs.Item($headingIndex + 3)\n\n    $p2.Range.Text = \"We create a copy of EMV2 properties in the EMV2 instance as owned property association of the respective EMV2 instance object. We do this for error events, error sources, error propagations.\"\n    $p3.Range.Text = \"Currently we create copies of Hazards, OccurrenceDistribution, Severity, and Likelihood properties. These are the only ones being used by safety analyses. There is a generic method to create an instance copy of any EMV2 property named as string.\"\n\n    # Insert two new paragraphs after p3: the first becomes the third new\n    # content paragraph, and the second remains empty, preserving the\n    # document's original trailing empty paragraph.\n    $p3.Range.InsertParagraphAfter()\n    $p4 = $d.Paragraphs.Item($headingIndex + 4)\n    $p4.Range.InsertParagraphAfter()\n    $p4.Range.Text = \"Note that property values can be associated with individual error types of the above mentioned EMV2 objects. We represent error events, error sources with multiple error types (listed as type sets) as separate ConstrainedInstanceObject for each type token representing each error type.  The property association is owned by the appropriate ConstrainedInstanceObject.\"\n}\n"}
